# Adds a new column BB to the table, carrying forward/refreshing the
# naive-forecaster combined YoY export with a new forecast vintage date
# (BB1 = 45986) and corresponding forecast values for rows 3-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date in BB1 - copy formatting/style from BA1 (style index 1)
# then overwrite with the new value.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Rows 3-18: new column BB just repeats the existing BA value for that row.
$sameValueRows = 3..18
foreach ($r in $sameValueRows) {
    $baCell = $ws.Cells.Item($r, 53)   # column BA = 53
    $bbCell = $ws.Cells.Item($r, 54)   # column BB = 54
    $bbCell.Value = $baCell.Value2
}

# Rows 19-21: new forecast values that differ from the BA column.
$ws.Cells.Item(19, 54).Value = -2.451276118722334
$ws.Cells.Item(20, 54).Value = -1.596682557877005
$ws.Cells.Item(21, 54).Value = -2.847551894053546
